$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70
$ws.Range("A70").Value = 131268322
$ws.Range("B70").Value = 83207
$ws.Range("D70").Value = 'LC'
$ws.Range("E70").Value = 6439
$ws.Range("F70").Value = 'Gulnål'
$ws.Range("G70").Value = 'Chaenotheca brachypoda'
$ws.Range("H70").Value = '(Ach.) Tibell'
$ws.Range("I70").NumberFormat = "@"
$ws.Range("I70").Value = ""
$ws.Range("P70").Value = 'Risboskogen, söder om, Upl'
$ws.Range("Q70").Value = 632054
$ws.Range("R70").Value = 6661554
$ws.Range("S70").Value = 15
$ws.Range("T70").Value = 'Uppsala'
$ws.Range("U70").Value = 'Uppsala'
$ws.Range("V70").Value = 'Uppland'
$ws.Range("W70").Value = 'Skuttunge'
$ws.Range("Y70").NumberFormat = "@"
$ws.Range("Y70").Value = '2026-02-22'
$ws.Range("AA70").NumberFormat = "@"
$ws.Range("AA70").Value = '2026-02-22'
$ws.Range("AD70").Value = $false
$ws.Range("AE70").Value = $false
$ws.Range("AG70").Value = $false
$ws.Range("AT70").NumberFormat = "@"
$ws.Range("AT70").Value = ""
$ws.Range("AW70").Value = 'Vilhelm Kroon'
$ws.Range("AX70").Value = 'Vilhelm Kroon'
$ws.Range("AY70").NumberFormat = "@"
$ws.Range("AY70").Value = ""

# Row 71
$ws.Range("A71").Value = 131268316
$ws.Range("B71").Value = 57881
$ws.Range("D71").Value = 'NT'
$ws.Range("E71").Value = 100049
$ws.Range("F71").Value = 'Spillkråka'
$ws.Range("G71").Value = 'Dryocopus martius'
$ws.Range("H71").Value = '(Linnaeus, 1758)'
$ws.Range("I71").NumberFormat = "@"
$ws.Range("I71").Value = ""
$ws.Range("K71").NumberFormat = "@"
$ws.Range("K71").Value = ""
$ws.Range("L71").NumberFormat = "@"
$ws.Range("L71").Value = ""
$ws.Range("M71").Value = 'äldre spår'
$ws.Range("N71").NumberFormat = "@"
$ws.Range("N71").Value = ""
$ws.Range("P71").Value = 'Risboskogen, söder om, Upl'
$ws.Range("Q71").Value = 632039
$ws.Range("R71").Value = 6661733
$ws.Range("S71").Value = 15
$ws.Range("T71").Value = 'Uppsala'
$ws.Range("U71").Value = 'Uppsala'
$ws.Range("V71").Value = 'Uppland'
$ws.Range("W71").Value = 'Skuttunge'
$ws.Range("Y71").NumberFormat = "@"
$ws.Range("Y71").Value = '2026-02-22'
$ws.Range("AA71").NumberFormat = "@"
$ws.Range("AA71").Value = '2026-02-22'
$ws.Range("AD71").Value = $false
$ws.Range("AE71").Value = $false
$ws.Range("AG71").Value = $false
$ws.Range("AT71").NumberFormat = "@"
$ws.Range("AT71").Value = ""
$ws.Range("AW71").Value = 'Vilhelm Kroon'
$ws.Range("AX71").Value = 'Vilhelm Kroon'
$ws.Range("AY71").NumberFormat = "@"
$ws.Range("AY71").Value = ""

# Row 72
$ws.Range("A72").Value = 131268315
$ws.Range("B72").Value = 57881
$ws.Range("D72").Value = 'NT'
$ws.Range("E72").Value = 100049
$ws.Range("F72").Value = 'Spillkråka'
$ws.Range("G72").Value = 'Dryocopus martius'
$ws.Range("H72").Value = '(Linnaeus, 1758)'
$ws.Range("I72").NumberFormat = "@"
$ws.Range("I72").Value = ""
$ws.Range("K72").NumberFormat = "@"
$ws.Range("K72").Value = ""
$ws.Range("L72").NumberFormat = "@"
$ws.Range("L72").Value = ""
$ws.Range("M72").Value = 'äldre spår'
$ws.Range("N72").NumberFormat = "@"
$ws.Range("N72").Value = ""
$ws.Range("P72").Value = 'Risboskogen, söder om, Upl'
$ws.Range("Q72").Value = 631848
$ws.Range("R72").Value = 6661810
$ws.Range("S72").Value = 15
$ws.Range("T72").Value = 'Uppsala'
$ws.Range("U72").Value = 'Uppsala'
$ws.Range("V72").Value = 'Uppland'
$ws.Range("W72").Value = 'Skuttunge'
$ws.Range("Y72").NumberFormat = "@"
$ws.Range("Y72").Value = '2026-02-22'
$ws.Range("AA72").NumberFormat = "@"
$ws.Range("AA72").Value = '2026-02-22'
$ws.Range("AD72").Value = $false
$ws.Range("AE72").Value = $false
$ws.Range("AG72").Value = $false
$ws.Range("AT72").NumberFormat = "@"
$ws.Range("AT72").Value = ""
$ws.Range("AW72").Value = 'Vilhelm Kroon'
$ws.Range("AX72").Value = 'Vilhelm Kroon'
$ws.Range("AY72").NumberFormat = "@"
$ws.Range("AY72").Value = ""

# Row 73
$ws.Range("A73").Value = 131268317
$ws.Range("B73").Value = 92268
$ws.Range("D73").Value = 'VU'
$ws.Range("E73").Value = 1209
$ws.Range("F73").Value = 'Rynkskinn'
$ws.Range("G73").Value = 'Hermanssonia centrifuga'
$ws.Range("H73").Value = '(P. Karst.) Zmitr.'
$ws.Range("I73").NumberFormat = "@"
$ws.Range("I73").Value = ""
$ws.Range("P73").Value = 'Risboskogen, söder om, Upl'
$ws.Range("Q73").Value = 632036
$ws.Range("R73").Value = 6661722
$ws.Range("S73").Value = 15
$ws.Range("T73").Value = 'Uppsala'
$ws.Range("U73").Value = 'Uppsala'
$ws.Range("V73").Value = 'Uppland'
$ws.Range("W73").Value = 'Skuttunge'
$ws.Range("Y73").NumberFormat = "@"
$ws.Range("Y73").Value = '2026-02-22'
$ws.Range("AA73").NumberFormat = "@"
$ws.Range("AA73").Value = '2026-02-22'
$ws.Range("AD73").Value = $false
$ws.Range("AE73").Value = $false
$ws.Range("AG73").Value = $false
$ws.Range("AT73").NumberFormat = "@"
$ws.Range("AT73").Value = ""
$ws.Range("AW73").Value = 'Vilhelm Kroon'
$ws.Range("AX73").Value = 'Vilhelm Kroon'
$ws.Range("AY73").NumberFormat = "@"
$ws.Range("AY73").Value = ""

# Row 74
$ws.Range("A74").Value = 131268314
$ws.Range("B74").Value = 97254
$ws.Range("D74").Value = 'NT'
$ws.Range("E74").Value = 53
$ws.Range("F74").Value = 'Vedtrappmossa'
$ws.Range("G74").Value = 'Crossocalyx hellerianus'
$ws.Range("H74").Value = '(Nees ex Lindenb.) Meyl.'
$ws.Range("I74").NumberFormat = "@"
$ws.Range("I74").Value = ""
$ws.Range("P74").Value = 'Risboskogen, söder om, Upl'
$ws.Range("Q74").Value = 632019
$ws.Range("R74").Value = 6661463
$ws.Range("S74").Value = 15
$ws.Range("T74").Value = 'Uppsala'
$ws.Range("U74").Value = 'Uppsala'
$ws.Range("V74").Value = 'Uppland'
$ws.Range("W74").Value = 'Skuttunge'
$ws.Range("Y74").NumberFormat = "@"
$ws.Range("Y74").Value = '2026-02-22'
$ws.Range("AA74").NumberFormat = "@"
$ws.Range("AA74").Value = '2026-02-22'
$ws.Range("AD74").Value = $false
$ws.Range("AE74").Value = $false
$ws.Range("AG74").Value = $false
$ws.Range("AT74").NumberFormat = "@"
$ws.Range("AT74").Value = ""
$ws.Range("AW74").Value = 'Vilhelm Kroon'
$ws.Range("AX74").Value = 'Vilhelm Kroon'
$ws.Range("AY74").NumberFormat = "@"
$ws.Range("AY74").Value = ""

# Row 75
$ws.Range("A75").Value = 131268321
$ws.Range("B75").Value = 57064
$ws.Range("D75").Value = 'NT'
$ws.Range("E75").Value = 102612
$ws.Range("F75").Value = 'Järpe'
$ws.Range("G75").Value = 'Tetrastes bonasia'
$ws.Range("H75").Value = '(Linnaeus, 1758)'
$ws.Range("I75").NumberFormat = "@"
$ws.Range("I75").Value = ""
$ws.Range("K75").NumberFormat = "@"
$ws.Range("K75").Value = ""
$ws.Range("L75").NumberFormat = "@"
$ws.Range("L75").Value = ""
$ws.Range("M75").Value = 'parning/parningsceremonier'
$ws.Range("N75").NumberFormat = "@"
$ws.Range("N75").Value = ""
$ws.Range("P75").Value = 'Risboskogen, söder om, Upl'
$ws.Range("Q75").Value = 632020
$ws.Range("R75").Value = 6661596
$ws.Range("S75").Value = 15
$ws.Range("T75").Value = 'Uppsala'
$ws.Range("U75").Value = 'Uppsala'
$ws.Range("V75").Value = 'Uppland'
$ws.Range("W75").Value = 'Skuttunge'
$ws.Range("Y75").NumberFormat = "@"
$ws.Range("Y75").Value = '2026-02-22'
$ws.Range("AA75").NumberFormat = "@"
$ws.Range("AA75").Value = '2026-02-22'
$ws.Range("AC75").Value = 'Framlockad med ljuduppspelning enl. vedertagen inventeringsmetodik för att konstatera revir.'
$ws.Range("AD75").Value = $false
$ws.Range("AE75").Value = $false
$ws.Range("AG75").Value = $false
$ws.Range("AH75").Value = 'Blandsumpskog'
$ws.Range("AI75").Value = 'Fuktig stråk med björk, klibbal och inväxt gran i naturskogsartad barrblandskog med lövinslag'
$ws.Range("AT75").NumberFormat = "@"
$ws.Range("AT75").Value = ""
$ws.Range("AW75").Value = 'Vilhelm Kroon'
$ws.Range("AX75").Value = 'Vilhelm Kroon'
$ws.Range("AY75").NumberFormat = "@"
$ws.Range("AY75").Value = ""
